$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17:C17").ClearContents()

$ws.Range("A17").Value = "regen_item_name"
$ws.Range("B17").Value = "Regeneracion."
$ws.Range("C17").Value = "Regen."

$ws.Range("A18").Value = "regen_item_desc"
$ws.Range("B18").Value = "Desc."
$ws.Range("C18").Value = "Desc."

$ws.Range("A19").Value = "regen_effect_desc"
$ws.Range("B19").Value = "Regenera {0} vida cada {1} segundos."
$ws.Range("C19").Value = "Recover {0} health every {1} seconds."

$ws.Range("C19").Select()
